$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 6.14
$ws.Range("D6").Value = -7.897
$ws.Range("B14").Value = 5.775
$ws.Range("C15").Value = -12.896
$ws.Range("B16").Value = 5.083
$ws.Range("D18").Value = -8.461
$ws.Range("D19").Value = -7.893000000000001
$ws.Range("B21").Value = 9.459000000000001
$ws.Range("C21").Value = -12.202
$ws.Range("C22").Value = -12.692
$ws.Range("B23").Value = 8.843999999999999
$ws.Range("C24").Value = -12.019
$ws.Range("B25").Value = 5.705
$ws.Range("B26").Value = 6.075
$ws.Range("C27").Value = -12.437
$ws.Range("C28").Value = -12.521
$ws.Range("B29").Value = 5.720999999999999
$ws.Range("D35").Value = -8.065999999999999
$ws.Range("C36").Value = -13.045
$ws.Range("C39").Value = -13.019
$ws.Range("B40").Value = 9.082000000000001
$ws.Range("D44").Value = -7.06
$ws.Range("C45").Value = -12.505
$ws.Range("D47").Value = -7.508
$ws.Range("C48").Value = -11.403
$ws.Range("C49").Value = -12.61
$ws.Range("D50").Value = -8.564
$ws.Range("D51").Value = -8.283000000000001
$ws.Range("C52").Value = -11.603
$ws.Range("D52").Value = -7.619
$ws.Range("B53").Value = 5.601000000000001
$ws.Range("C53").Value = -10.446
$ws.Range("C54").Value = -12.418
$ws.Range("D55").Value = -8.499000000000001
$ws.Range("B57").Value = 5.673
$ws.Range("C57").Value = -13.019
$ws.Range("D57").Value = -8.456999999999999
$ws.Range("D58").Value = -8.175000000000001
$ws.Range("B59").Value = 5.274
$ws.Range("D64").Value = -7.816999999999998
$ws.Range("B65").Value = 6.068
$ws.Range("D66").Value = -7.627
$ws.Range("B69").Value = 6.233
$ws.Range("C70").Value = -11.581
$ws.Range("C71").Value = -10.85
$ws.Range("B79").Value = 6.759
$ws.Range("D80").Value = -7.872999999999999
$ws.Range("B83").Value = 5.359
$ws.Range("D83").Value = -8.735000000000001
$ws.Range("C86").Value = -13.339
$ws.Range("C87").Value = -13.474
$ws.Range("C89").Value = -13.153
$ws.Range("B91").Value = 6.14
$ws.Range("D92").Value = -7.058
$ws.Range("B93").Value = 6.303
$ws.Range("D94").Value = -7.184
$ws.Range("D96").Value = -7.35
$ws.Range("D97").Value = -7.290999999999999
$ws.Range("B100").Value = 5.525
$ws.Range("C101").Value = -12.083
$ws.Range("D101").Value = -7.725
$ws.Range("B103").Value = 5.684
